# Apply the "Regional Availability Factor" update to the RAF-generation sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RAF-generation")

# --- Update capacity-factor values in column B ---
$ws.Range("B2").Value  = 0.46666666699999998   # hard coal
$ws.Range("B3").Value  = 0.7                   # natural gas steam turbine
$ws.Range("B4").Value  = 0.7                   # natural gas combustion turbine
$ws.Range("B12").Value = 0.7                   # petroleum
$ws.Range("B13").Value = 0.7                   # natural gas peaker
$ws.Range("B15").Value = 1                     # offshore wind
$ws.Range("B16").Value = 0.7                   # crude oil
$ws.Range("B17").Value = 0.7                   # heavy or residual fuel oil
$ws.Range("B18").Value = 0.7                   # municipal solid waste

# --- Header cell A1: drop the italic/word-wrap style and let row 1 resume its default height ---
$ws.Cells.Item(1, 1).Style = "Normal"
$ws.Rows.Item(1).AutoFit()

# --- Update the selection shown when the workbook is reopened ---
$ws.Range("A1:B25").Select()
